$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderInfo")

# Sales Tax: 10% -> 7%
$ws.Range("K2").Value = 0.07

# Delivery Day: 22 -> 21
$ws.Range("Q2").Value = 21

# Olyve Premiere Code: "nope" -> "beauty10"
$ws.Range("Y2").Value = "beauty10"
